$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the stale external workbook link (no longer referenced) ---
$links = $wb.LinkSources(1)
if ($links) {
    foreach ($link in $links) {
        $wb.BreakLink($link, 1)
    }
}

# --- Header / title area updates ---
# A2 used to read "SW DIPLOMA" -> now "Proyecto Integrador"
$ws.Range("A2").Value = "Proyecto Integrador"
# A3 used to read "Project Lead" -> now blank
$ws.Range("A3").Value = ""

# --- % DONE column updates for later rows ---
$ws.Range("G20").Value = 80
$ws.Range("G22").Value = 85
$ws.Range("G23").Value = 90
$ws.Range("G24").Value = 100

# --- View settings: zoom + top-left cell + selection ---
$ws.Application.ActiveWindow.Zoom = 60
$ws.Range("G21").Select()
